$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H4").Value = 100.25
$ws.Range("J4").Value = 41
$ws.Range("L4").Value = 41
$ws.Range("N4").Value = -269
$ws.Range("H74").Value = 8084.143
$ws.Range("I74").Value = 8084.143
$ws.Range("K74").Value = 8084.143
$ws.Range("M74").Value = -7148.143
$ws.Range("H77").Value = 8084.143
$ws.Range("I77").Value = 8084.143
$ws.Range("K77").Value = 40420.715
$ws.Range("M77").Value = -35740.715
$ws.Range("H93").Value = 99999
$ws.Range("J93").Value = 99999
$ws.Range("L93").Value = 99999
$ws.Range("N93").Value = -104991
$ws.Range("H100").Value = 37824.57
$ws.Range("I100").Value = 42053.56
$ws.Range("J100").Value = 2583
$ws.Range("K100").Value = 42053.56
$ws.Range("L100").Value = 2583
$ws.Range("M100").Value = -41512.56
$ws.Range("N100").Value = -3665
$ws.Range("H137").Value = 11548.069
$ws.Range("I137").Value = 5316.737
$ws.Range("J137").Value = 16481.209
$ws.Range("K137").Value = 15950.211
$ws.Range("L137").Value = 49443.62699999999
$ws.Range("M137").Value = -13400.211
$ws.Range("N137").Value = -54543.62699999999
$ws.Range("H141").Value = 3398.8572
$ws.Range("I141").Value = 3447.8333
$ws.Range("K141").Value = 10343.4999
$ws.Range("M141").Value = -5163.499899999999

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H60").Value = 35992.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 35992.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 35992.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -37458.5
$ws.Range("H61").Value = 5543.933
$ws.Range("I61").Value = 3205.5454
$ws.Range("K61").Value = 3205.5454
$ws.Range("M61").Value = -2993.5454
$ws.Range("H74").Value = 4116.4165
$ws.Range("J74").Value = 6859.4
$ws.Range("L74").Value = 6859.4
$ws.Range("N74").Value = -8607.4
$ws.Range("H77").Value = 4116.4165
$ws.Range("J77").Value = 6859.4
$ws.Range("L77").Value = 34297
$ws.Range("N77").Value = -43033
$ws.Range("H132").Value = 6600.0264
$ws.Range("I132").Value = 6866.028
$ws.Range("K132").Value = 20598.084
$ws.Range("M132").Value = -18068.084
$ws.Range("H136").Value = 5543.933
$ws.Range("I136").Value = 3205.5454
$ws.Range("K136").Value = 9616.636200000001
$ws.Range("M136").Value = -7066.636200000001

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 1468.8948
$ws.Range("J107").Value = 1871
$ws.Range("L107").Value = 1871
$ws.Range("N107").Value = -5711
$ws.Range("H134").Value = 8777
$ws.Range("I134").Value = 6570.0625
$ws.Range("J134").Value = 12308.1
$ws.Range("K134").Value = 19710.1875
$ws.Range("L134").Value = 36924.3
$ws.Range("M134").Value = -17175.1875
$ws.Range("N134").Value = -41994.3

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2397.318
$ws.Range("I31").Value = 1685.3529
$ws.Range("J31").Value = 4818
$ws.Range("K31").Value = 1685.3529
$ws.Range("L31").Value = 4818
$ws.Range("M31").Value = -1390.3529
$ws.Range("N31").Value = -5408
$ws.Range("H34").Value = 2397.318
$ws.Range("I34").Value = 1685.3529
$ws.Range("J34").Value = 4818
$ws.Range("K34").Value = 1685.3529
$ws.Range("L34").Value = 4818
$ws.Range("M34").Value = -1483.3529
$ws.Range("N34").Value = -5222
$ws.Range("H132").Value = 24911.318
$ws.Range("I132").Value = 15322
$ws.Range("K132").Value = 45966
$ws.Range("M132").Value = -43436

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H75").Value = 66668492
$ws.Range("J75").Value = 83335090
$ws.Range("L75").Value = 250005270
$ws.Range("N75").Value = -250007266
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H78").Value = 66668492
$ws.Range("J78").Value = 83335090
$ws.Range("L78").Value = 750015810
$ws.Range("N78").Value = -750025794
$ws.Range("H80").Value = 2200
$ws.Range("J80").Value = 2200
$ws.Range("L80").Value = 6600
$ws.Range("N80").Value = -8472
$ws.Range("H81").Value = 35329.465
$ws.Range("I81").Value = 42916.5
$ws.Range("J81").Value = 4981.3335
$ws.Range("K81").Value = 128749.5
$ws.Range("L81").Value = 14944.0005
$ws.Range("M81").Value = -127626.5
$ws.Range("N81").Value = -17190.0005
$ws.Range("H83").Value = 2200
$ws.Range("J83").Value = 2200
$ws.Range("L83").Value = 19800
$ws.Range("N83").Value = -29160
$ws.Range("H84").Value = 35329.465
$ws.Range("I84").Value = 42916.5
$ws.Range("J84").Value = 4981.3335
$ws.Range("K84").Value = 386248.5
$ws.Range("L84").Value = 44832.0015
$ws.Range("M84").Value = -380632.5
$ws.Range("N84").Value = -56064.0015
$ws.Range("H87").Value = 8675.333000000001
$ws.Range("I87").Value = 8675.333000000001
$ws.Range("K87").Value = 26025.999
$ws.Range("M87").Value = -24777.999
$ws.Range("H88").Value = 10055.3125
$ws.Range("J88").Value = 10028.064
$ws.Range("L88").Value = 30084.192
$ws.Range("N88").Value = -30940.192
$ws.Range("H90").Value = 8675.333000000001
$ws.Range("I90").Value = 8675.333000000001
$ws.Range("K90").Value = 78077.997
$ws.Range("M90").Value = -71837.997
$ws.Range("H91").Value = 10055.3125
$ws.Range("J91").Value = 10028.064
$ws.Range("L91").Value = 30084.192
$ws.Range("N91").Value = -33048.192

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 3028.1667
$ws.Range("I132").Value = 2599.9285
$ws.Range("J132").Value = 4527
$ws.Range("K132").Value = 7799.7855
$ws.Range("L132").Value = 13581
$ws.Range("M132").Value = -5269.7855
$ws.Range("N132").Value = -18641

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 526.13336
$ws.Range("I55").Value = 1246
$ws.Range("J55").Value = 264.36365
$ws.Range("K55").Value = 1246
$ws.Range("L55").Value = 264.36365
$ws.Range("M55").Value = -1073
$ws.Range("N55").Value = -610.36365
$ws.Range("H68").Value = 2838.6428
$ws.Range("J68").Value = 5664.6665
$ws.Range("L68").Value = 5664.6665
$ws.Range("N68").Value = -7162.6665
$ws.Range("H71").Value = 2838.6428
$ws.Range("J71").Value = 5664.6665
$ws.Range("L71").Value = 28323.3325
$ws.Range("N71").Value = -35811.3325
$ws.Range("H122").Value = 5050.5713
$ws.Range("I122").Value = 4893.3335
$ws.Range("J122").Value = 5994
$ws.Range("K122").Value = 14680.0005
$ws.Range("L122").Value = 17982
$ws.Range("M122").Value = -12230.0005
$ws.Range("N122").Value = -22882
$ws.Range("H136").Value = 5850410.5
$ws.Range("I136").Value = 7409661
$ws.Range("K136").Value = 22228983
$ws.Range("M136").Value = -22226433

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 490.6316
$ws.Range("I113").Value = 390.33334
$ws.Range("J113").Value = 662.5714
$ws.Range("K113").Value = 1171.00002
$ws.Range("L113").Value = 1987.7142
$ws.Range("M113").Value = 998.9999800000001
$ws.Range("N113").Value = -6327.7142
$ws.Range("H122").Value = 2945.8215
$ws.Range("I122").Value = 2977.1853
$ws.Range("K122").Value = 8931.555899999999
$ws.Range("M122").Value = -6481.555899999999
$ws.Range("H132").Value = 30505.033
$ws.Range("I132").Value = 25234.7
$ws.Range("J132").Value = 40087.453
$ws.Range("K132").Value = 75704.10000000001
$ws.Range("L132").Value = 120262.359
$ws.Range("M132").Value = -73174.10000000001
$ws.Range("N132").Value = -125322.359
$ws.Range("H136").Value = 1314.6154
$ws.Range("J136").Value = 4999.5
$ws.Range("L136").Value = 14998.5
$ws.Range("N136").Value = -20098.5
